$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "done" for several tasks that already have an assignee, in the new column C.
$ws.Range("C2").Value = "done"
$ws.Range("C3").Value = "done"
$ws.Range("C6").Value = "done"
$ws.Range("C8").Value = "done"
$ws.Range("C9").Value = "done"

# Fix assignee for "Missions 2.1: limite inventaire" row (row 7) from Roméo to Fabio
$ws.Range("B7").Value = "Fabio"

# A6 had a redundant duplicated style (applyFont + center alignment) - normalize it to
# just the plain centered style used elsewhere, by re-asserting the center alignment.
$ws.Range("A6").HorizontalAlignment = -4108

# Update selection to C1
$ws.Range("C1").Select()
